# "Edit station linked with uploaded excel"
# The schedule's timestamp (A2) is refreshed to the newest upload time,
# and the current selection moves to the freshly-updated cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the linked timestamp to the new upload time (12/7/2022 3:17:00 PM)
$ws.Range("A2").Value = 44902.636805555558

# Reflect the edited cell as the active selection
$ws.Range("A2").Select()
